# Refresh cryptos list - updates Price (D) and Volume/1h (E) text columns for
# each coin row, plus a two-row swap of FraxShare/Algorand (rows 36-37).
# All D/E cells hold plain text (e.g. "28.494.50", "  +0.66%  "), so each
# write forces NumberFormat "@" first to stop values that look numeric
# (e.g. "0.6370", "8.903") from being auto-coerced into numbers, then
# restores the cell's original (unstyled) Style so no stray style index is
# left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with default (unstyled) format, used to reset style after
# forcing text NumberFormat -- keeps cells free of stray quote-prefix / text styles
$styleRef = $ws.Range("A1")

function Set-TextValue($ref, $value) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value2 = $value
    $cell.Style = $styleRef.Style
}

Set-TextValue "D2" '28.494.50'
Set-TextValue "E2" '  +0.66%  '
Set-TextValue "D3" '1.873.62'
Set-TextValue "E3" '  +0.60%  '
Set-TextValue "E4" '  -1.02%  '
Set-TextValue "D5" '315.45'
Set-TextValue "E5" '  +0.19%  '
Set-TextValue "E6" '  -0.71%  '
Set-TextValue "D7" '0.5092'
Set-TextValue "E7" '  -0.23%  '
Set-TextValue "D8" '0.3897'
Set-TextValue "E8" '  -0.60%  '
Set-TextValue "D9" '0.08342'
Set-TextValue "E9" '  +0.74%  '
Set-TextValue "D10" '1.104'
Set-TextValue "E10" '  -0.55%  '
Set-TextValue "D11" '6.224'
Set-TextValue "E11" '  +0.06%  '
Set-TextValue "D12" '1.878.50'
Set-TextValue "E12" '  -0.26%  '
Set-TextValue "D13" '20.39'
Set-TextValue "E13" '  +0.49%  '
Set-TextValue "D14" '7.271'
Set-TextValue "E14" '  +0.70%  '
Set-TextValue "D15" '1.008'
Set-TextValue "E15" '  -0.81%  '
Set-TextValue "E16" '  +0.22%  '
Set-TextValue "D17" '91.20'
Set-TextValue "E17" '  -0.04%  '
Set-TextValue "E18" '  +0.22%  '
Set-TextValue "D19" '17.67'
Set-TextValue "E19" '  +0.67%  '
Set-TextValue "E20" '  -0.81%  '
Set-TextValue "D21" '5.905'
Set-TextValue "E21" '  -0.45%  '
Set-TextValue "D22" '28.522.33'
Set-TextValue "E22" '  +0.68%  '
Set-TextValue "E23" '  +0.11%  '
Set-TextValue "D24" '2.225'
Set-TextValue "E24" '  -1.26%  '
Set-TextValue "D25" '2.091.58'
Set-TextValue "E25" '  -0.12%  '
Set-TextValue "D26" '161.31'
Set-TextValue "E26" '  +0.63%  '
Set-TextValue "D27" '20.60'
Set-TextValue "E27" '  -0.32%  '
Set-TextValue "D28" '2.428'
Set-TextValue "E28" '  +3.01%  '
Set-TextValue "D29" '127.12'
Set-TextValue "E29" '  +0.41%  '
Set-TextValue "E30" '  -0.17%  '
Set-TextValue "D31" '1.038'
Set-TextValue "E31" '  +1.34%  '
Set-TextValue "D32" '5.734'
Set-TextValue "E32" '  -1.17%  '
Set-TextValue "D33" '3.598'
Set-TextValue "E33" '  -1.00%  '
Set-TextValue "D34" '0.02455'
Set-TextValue "E34" '  +1.44%  '
Set-TextValue "D35" '0.06558'
Set-TextValue "E35" '  +1.85%  '
Set-TextValue "B36" 'FraxShare'
Set-TextValue "C36" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D36" '8.903'
Set-TextValue "E36" '  -1.94%  '
Set-TextValue "B37" 'Algorand'
Set-TextValue "C37" 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue "D37" '0.2155'
Set-TextValue "E37" '  -0.53%  '
Set-TextValue "D38" '5.032'
Set-TextValue "E38" '  +1.95%  '
Set-TextValue "E39" '  +0.42%  '
Set-TextValue "D40" '1.238'
Set-TextValue "E40" '  -0.33%  '
Set-TextValue "D41" '0.6370'
Set-TextValue "E41" '  -0.40%  '
Set-TextValue "E42" '  -0.03%  '
Set-TextValue "E43" '  -0.80%  '
Set-TextValue "D44" '0.5989'
Set-TextValue "E44" '  -0.14%  '
Set-TextValue "E45" '  +1.50%  '
Set-TextValue "E46" '  -0.03%  '
Set-TextValue "D47" '1.999'
Set-TextValue "E47" '  +1.20%  '
Set-TextValue "D48" '1.217'
Set-TextValue "D49" '122.00'
Set-TextValue "E49" '  +0.67%  '
Set-TextValue "D50" '1.147'
Set-TextValue "E50" '  -10.28%  '
Set-TextValue "D51" '0.06828'
Set-TextValue "E51" '  +0.00%  '
